$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (hunk 0)
$ws.Range("H2").Value = 947.8570999999999
$ws.Range("I2").Value = 913.625
$ws.Range("J2").Value = 993.5
$ws.Range("K2").Value = 913.625
$ws.Range("L2").Value = 993.5
$ws.Range("M2").Value = -800.625
$ws.Range("N2").Value = -1219.5

# Row 103 (hunk 1)
$ws.Range("H103").Value = 771.4286
$ws.Range("I103").Value = 1200
$ws.Range("J103").Value = 700
$ws.Range("K103").Value = 3600
$ws.Range("L103").Value = 2100
$ws.Range("M103").Value = -3014
$ws.Range("N103").Value = -3272

# Row 137 (hunk 2)
$ws.Range("H137").Value = 1423
$ws.Range("I137").Value = 1269.1333
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 3807.3999
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -1257.3999

# Row 141 (hunk 3)
$ws.Range("H141").Value = 2556
$ws.Range("I141").Value = 1560.52
$ws.Range("J141").Value = 14999.5
$ws.Range("K141").Value = 4681.559999999999
$ws.Range("L141").Value = 44998.5
$ws.Range("M141").Value = 498.4400000000005


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45 (hunk 4)
$ws.Range("H45").Value = 2666.5454
$ws.Range("I45").Value = 2683.2
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 2683.2
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -2306.2

# Row 74 (hunk 5)
$ws.Range("H74").Value = 1101.091
$ws.Range("I74").Value = 1061.2
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1061.2
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -187.2

# Row 77 (hunk 6)
$ws.Range("H77").Value = 1101.091
$ws.Range("I77").Value = 1061.2
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 5306
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -938

# Row 124 (hunk 7)
$ws.Range("H124").Value = 37500
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 37500
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 37500
$ws.Range("N124").Value = -47320

# Row 132 (hunk 8)
$ws.Range("H132").Value = 1183.6522
$ws.Range("I132").Value = 1210.5
$ws.Range("J132").Value = 1087
$ws.Range("K132").Value = 3631.5
$ws.Range("L132").Value = 3261
$ws.Range("M132").Value = -1101.5


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20 (hunk 9)
$ws.Range("H20").Value = 3959.8
$ws.Range("I20").Value = 5266.3335
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 5266.3335
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = -5019.3335

# Row 24 (hunk 10)
$ws.Range("H24").Value = 3599.5
$ws.Range("I24").Value = 3599.5
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 3599.5
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -3364.5
$ws.Range("N24").ClearContents()

# Row 86 (hunk 11)
$ws.Range("H86").Value = 1123.75
$ws.Range("I86").Value = 1165
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 1165
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -42
$ws.Range("N86").Value = -3246

# Row 89 (hunk 12)
$ws.Range("H89").Value = 1123.75
$ws.Range("I89").Value = 1165
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 5825
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -209
$ws.Range("N89").Value = -16232

# Row 105 (hunk 13)
$ws.Range("H105").Value = 2555.8333
$ws.Range("I105").Value = 2442.7693
$ws.Range("J105").Value = 2849.8
$ws.Range("K105").Value = 2442.7693
$ws.Range("L105").Value = 2849.8
$ws.Range("M105").Value = -695.7692999999999


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (hunk 14)
$ws.Range("H31").Value = 2853.1538
$ws.Range("I31").Value = 2044.65
$ws.Range("J31").Value = 5548.1665
$ws.Range("K31").Value = 2044.65
$ws.Range("L31").Value = 5548.1665
$ws.Range("M31").Value = -1749.65
$ws.Range("N31").Value = -6138.1665

# Row 34 (hunk 15)
$ws.Range("H34").Value = 2853.1538
$ws.Range("I34").Value = 2044.65
$ws.Range("J34").Value = 5548.1665
$ws.Range("K34").Value = 2044.65
$ws.Range("L34").Value = 5548.1665
$ws.Range("M34").Value = -1842.65
$ws.Range("N34").Value = -5952.1665

# Row 41 (hunk 16)
$ws.Range("H41").Value = 27091.666
$ws.Range("I41").Value = 1500
$ws.Range("J41").Value = 29418.182
$ws.Range("K41").Value = 1500
$ws.Range("L41").Value = 29418.182
$ws.Range("M41").Value = -1072

# Row 107 (hunk 17)
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()

# Row 132 (hunk 18)
$ws.Range("H132").Value = 3199.5
$ws.Range("I132").Value = 3199.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9598.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7068.5

# Row 134 (hunk 19)
$ws.Range("H134").Value = 3982.1428

# Row 138 (hunk 20)
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 21 (hunk 21)
$ws.Range("H21").Value = 5000
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -4827

# Row 30 (hunk 22)
$ws.Range("H30").Value = 5000
$ws.Range("I30").Value = 5000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 5000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -4895

# Row 70 (hunk 23)
$ws.Range("H70").Value = 5916.3335
$ws.Range("I70").Value = 5125.875
$ws.Range("J70").Value = 7497.25
$ws.Range("K70").Value = 5125.875
$ws.Range("L70").Value = 7497.25
$ws.Range("M70").Value = -4855.875

# Row 73 (hunk 24)
$ws.Range("H73").Value = 5916.3335
$ws.Range("I73").Value = 5125.875
$ws.Range("J73").Value = 7497.25
$ws.Range("K73").Value = 5125.875
$ws.Range("L73").Value = 7497.25
$ws.Range("M73").Value = -4189.875

# Row 132 (hunk 25)
$ws.Range("H132").Value = 2251.182
$ws.Range("I132").Value = 1990.5714
$ws.Range("J132").Value = 2707.25
$ws.Range("K132").Value = 5971.7142
$ws.Range("L132").Value = 8121.75
$ws.Range("M132").Value = -3441.7142
$ws.Range("N132").Value = -13181.75


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61 (hunk 26)
$ws.Range("H61").Value = 4371.8184
$ws.Range("I61").Value = 4309
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 4309
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -4107

# Row 64 (hunk 27)
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67 (hunk 28)
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 100 (hunk 29)
$ws.Range("H100").Value = 1225.75
$ws.Range("I100").Value = 967.6667
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 967.6667
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -426.6667
$ws.Range("N100").Value = -3082

# Row 113 (hunk 30)
$ws.Range("H113").Value = 4371.8184
$ws.Range("I113").Value = 4309
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 4309
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2139

# Row 122 (hunk 31)
$ws.Range("H122").Value = 24998.5
$ws.Range("I122").Value = 24998.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 74995.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -72545.5

# Row 132 (hunk 32)
$ws.Range("H132").Value = 2352.258
$ws.Range("I132").Value = 2035.7307
$ws.Range("J132").Value = 3998.2
$ws.Range("K132").Value = 6107.1921
$ws.Range("L132").Value = 11994.6
$ws.Range("M132").Value = -3577.1921

# Row 136 (hunk 33)
$ws.Range("H136").Value = 6821.2
$ws.Range("I136").Value = 6428.25
$ws.Range("J136").Value = 8393
$ws.Range("K136").Value = 19284.75
$ws.Range("L136").Value = 25179
$ws.Range("M136").Value = -16734.75


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132 (hunk 34)
$ws.Range("H132").Value = 3050.516
$ws.Range("I132").Value = 2766.9092
$ws.Range("J132").Value = 3743.7778
$ws.Range("K132").Value = 8300.7276
$ws.Range("L132").Value = 11231.3334
$ws.Range("M132").Value = -5770.7276

